# Adds a new "2022-Q4" sheet (fund holdings detail) before the existing
# "2022-Q3" sheet, and updates the "总计" (summary) sheet with the new
# 2022-Q4 row, pushing the older rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet (sheet index 1) with the new data.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryRows = @(
    @(0, "2022-Q4", 14, 4.83),
    @(1, "2022-Q3", 2,  0),
    @(2, "2022-Q2", 9,  1.79),
    @(3, "2022-Q1", 54, 8.42),
    @(4, "2021-Q4", 14, 1.66),
    @(5, "2021-Q3", 32, 5.52),
    @(6, "2021-Q2", 3,  0.93),
    @(7, "2020-Q4", 1,  0)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right before "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$fundRows = @(
    @(0,  "166301", "华商新趋势优选灵活配置混合",        "98.72", "74.44", "2.05", "2.0238", 3),
    @(1,  "000390", "华商优势行业混合",                  "35.74", "90.63", "2.21", "0.7899", 7),
    @(2,  "010977", "华夏鸿阳6个月持有期混合A",           "25.61", "84.81", "3.02", "0.7734", 5),
    @(3,  "630002", "华商盛世成长混合",                  "23.51", "90.49", "3.18", "0.7476", 5),
    @(4,  "003624", "创金合信资源主题精选股票A",          "3.79",  "91.89", "3.28", "0.1243", 9),
    @(5,  "010978", "华夏鸿阳6个月持有期混合C",           "3.25",  "84.81", "3.02", "0.0982", 5),
    @(6,  "003625", "创金合信资源主题精选股票C",          "2.82",  "91.89", "3.28", "0.0925", 9),
    @(7,  "008488", "华商恒益稳健混合",                  "4.22",  "49.85", "1.54", "0.0650", 6),
    @(8,  "481017", "工银量化策略混合A",                 "2.42",  "92.03", "2.08", "0.0503", 8),
    @(9,  "160620", "鹏华中证A股资源产业指数（LOF）A",     "1.72",  "94.49", "2.28", "0.0392", 1),
    @(10, "011888", "民生加银周期优选混合型证券投资基金A", "0.35",  "92.72", "3.52", "0.0123", 7),
    @(11, "012808", "鹏华中证A股资源产业指数（LOF）C",     "0.24",  "94.49", "2.28", "0.0055", 1),
    @(12, "011889", "民生加银周期优选混合型证券投资基金C", "0.09",  "92.72", "3.52", "0.0032", 7),
    @(13, "012241", "工银量化策略混合C",                 "0.01",  "92.03", "2.08", "0.0002", 8)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    # Column B (fund code) is text with leading zeros in the source
    # workbook (e.g. "000390") - force text so COM doesn't coerce it to
    # a number and drop the leading zeros.
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    # Columns D..G hold text that looks numeric in the source workbook
    # (stored as inlineStr, not numbers) - force text with a leading
    # apostrophe so COM doesn't silently coerce them to doubles.
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

Write-Output "2022-Q4 sheet added and 总计 sheet updated"
